$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H86").Value = 863.9375
$ws.Range("I86").Value = 898.5
$ws.Range("J86").Value = 806.3333
$ws.Range("K86").Value = 898.5
$ws.Range("L86").Value = 806.3333
$ws.Range("M86").Value = 224.5
$ws.Range("N86").Value = -3052.3333
$ws.Range("H89").Value = 863.9375
$ws.Range("I89").Value = 898.5
$ws.Range("J89").Value = 806.3333
$ws.Range("K89").Value = 4492.5
$ws.Range("L89").Value = 4031.6665
$ws.Range("M89").Value = 1123.5
$ws.Range("N89").Value = -15263.6665
$ws.Range("H106").Value = 1706.125
$ws.Range("I106").Value = 1708.1666
$ws.Range("K106").Value = 1708.1666
$ws.Range("M106").Value = -1077.1666
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H132").Value = 990.5
$ws.Range("I132").Value = 989.55
$ws.Range("K132").Value = 2968.65
$ws.Range("M132").Value = -438.6499999999996
$ws.Range("H138").Value = 3211.0715
$ws.Range("J138").Value = 2118.2778
$ws.Range("L138").Value = 6354.8334
$ws.Range("N138").Value = -16634.8334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2263.2727
$ws.Range("I45").Value = 1563.5
$ws.Range("K45").Value = 1563.5
$ws.Range("M45").Value = -1186.5
$ws.Range("H57").Value = 8000
$ws.Range("I57").Value = 8000
$ws.Range("K57").Value = 8000
$ws.Range("M57").Value = -7516
$ws.Range("H97").Value = 3010
$ws.Range("I97").Value = 3010
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3010
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2514
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 4099.769
$ws.Range("I132").Value = 3820.5557
$ws.Range("K132").Value = 11461.6671
$ws.Range("M132").Value = -8931.667099999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 266.66666
$ws.Range("I22").Value = 266.66666
$ws.Range("K22").Value = 266.66666
$ws.Range("M22").Value = -93.66665999999998
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H105").Value = 2438.5652
$ws.Range("I105").Value = 2438.5652
$ws.Range("K105").Value = 2438.5652
$ws.Range("M105").Value = -691.5652
$ws.Range("H107").Value = 3049.0625
$ws.Range("I107").Value = 3108.7273
$ws.Range("K107").Value = 3108.7273
$ws.Range("M107").Value = -1188.7273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1629.2174
$ws.Range("I31").Value = 1052.8182
$ws.Range("J31").Value = 2157.5833
$ws.Range("K31").Value = 1052.8182
$ws.Range("L31").Value = 2157.5833
$ws.Range("M31").Value = -757.8181999999999
$ws.Range("N31").Value = -2747.5833
$ws.Range("H34").Value = 1629.2174
$ws.Range("I34").Value = 1052.8182
$ws.Range("J34").Value = 2157.5833
$ws.Range("K34").Value = 1052.8182
$ws.Range("L34").Value = 2157.5833
$ws.Range("M34").Value = -850.8181999999999
$ws.Range("N34").Value = -2561.5833
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H116").Value = 100001976
$ws.Range("I116").Value = 1328
$ws.Range("J116").Value = 125002136
$ws.Range("K116").Value = 3984
$ws.Range("L116").Value = 375006408
$ws.Range("M116").Value = -542
$ws.Range("N116").Value = -375013292
$ws.Range("H131").Value = 16691332
$ws.Range("J131").Value = 36795.95
$ws.Range("L131").Value = 110387.85
$ws.Range("N131").Value = -120467.85

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2060531.6
$ws.Range("I126").Value = 3270853.8
$ws.Range("J126").Value = 2984.1
$ws.Range("K126").Value = 9812561.399999999
$ws.Range("L126").Value = 8952.299999999999
$ws.Range("M126").Value = -9810091.399999999
$ws.Range("N126").Value = -13892.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4300.4287
$ws.Range("I7").Value = 2601.1428
$ws.Range("J7").Value = 5999.7144
$ws.Range("K7").Value = 2601.1428
$ws.Range("L7").Value = 5999.7144
$ws.Range("M7").Value = -2489.1428
$ws.Range("N7").Value = -6223.7144
$ws.Range("H22").Value = 2059.9
$ws.Range("I22").Value = 1433.1666
$ws.Range("K22").Value = 1433.1666
$ws.Range("M22").Value = -1138.1666
$ws.Range("H25").Value = 9453.5
$ws.Range("I25").Value = 9407
$ws.Range("K25").Value = 9407
$ws.Range("M25").Value = -9177
$ws.Range("H27").Value = 2059.9
$ws.Range("I27").Value = 1433.1666
$ws.Range("K27").Value = 1433.1666
$ws.Range("M27").Value = -1326.1666
$ws.Range("H55").Value = 121.35714
$ws.Range("I55").Value = 132.83333
$ws.Range("K55").Value = 132.83333
$ws.Range("M55").Value = 40.16667000000001
$ws.Range("H126").Value = 4300.4287
$ws.Range("I126").Value = 2601.1428
$ws.Range("J126").Value = 5999.7144
$ws.Range("K126").Value = 7803.428400000001
$ws.Range("L126").Value = 17999.1432
$ws.Range("M126").Value = -5333.428400000001
$ws.Range("N126").Value = -22939.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 10000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 10000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10000
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -10554
$ws.Range("H107").Value = 1771.625
$ws.Range("I107").Value = 1797.6666
$ws.Range("K107").Value = 5392.9998
$ws.Range("L107").Value = 5268
$ws.Range("M107").Value = -3472.9998
$ws.Range("N107").Value = -9108
$ws.Range("H132").Value = 1176.8
$ws.Range("I132").Value = 502.82352
$ws.Range("K132").Value = 1508.47056
$ws.Range("M132").Value = 1021.52944
$ws.Range("H135").Value = 99818.836
$ws.Range("J135").Value = 99818.836
$ws.Range("L135").Value = 99818.836
$ws.Range("N135").Value = -109958.836
$ws.Range("H136").Value = 34726250
$ws.Range("I136").Value = 61732612
$ws.Range("K136").Value = 185197836
$ws.Range("M136").Value = -185195286

Write-Output "done"